# Rename existing "Experiment_2" sheet to "Experiment_1" and update its data,
# then add a new "Experiment_4" sheet (placed after the last sheet) with its own data.

$wb = $excel.ActiveWorkbook

# --- Update existing sheet "Experiment_2" -> rename + new values ---
$ws2 = $wb.Worksheets.Item("Experiment_2")
$ws2.Name = "Experiment_1"
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 17.90007758140564

# --- Add a new sheet "Experiment_4" after the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Experiment_4"

$newSheet.Range("A1").Value = "Similarity"
$newSheet.Range("B1").Value = "Inference_Time"

# Match the header formatting used on the other experiment sheets.
$ws2.Range("A1:B1").Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)

$newSheet.Range("A2").Value = 6
$newSheet.Range("B2").Value = 193.8237497806549
